# Fruta / hortaliza, semanal
#
# A new weekly observation was inserted for "Feria Lagunitas de Puerto Montt -
# Pepino dulce" at row 55 (pushing the existing rows 55-74 down to 56-75,
# which simply shifts each pre-existing record one row down and re-appends
# the old last row (74) as the new last row (75)). The brand new record gets
# its own data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 55..74 down to 56..75 and open up a blank row 55 for the new
# weekly record (mirrors Excel's "Insert Copied Cells"/row insert behaviour,
# carrying the row-above formatting, e.g. the date style on column D).
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly record.
$ws.Cells.Item(55, 1).Value = 4
$ws.Cells.Item(55, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(55, 3).Value = "Los Lagos"
$ws.Cells.Item(55, 4).Value = 44988
$ws.Cells.Item(55, 5).Value = 10
$ws.Cells.Item(55, 6).Value = 100112043
$ws.Cells.Item(55, 7).Value = "Pepino dulce"
$ws.Cells.Item(55, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 100
$ws.Cells.Item(55, 11).Value = 20000
$ws.Cells.Item(55, 12).Value = 22000
$ws.Cells.Item(55, 13).Value = 21000
$ws.Cells.Item(55, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(55, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(55, 16).Value = 1167
$ws.Cells.Item(55, 17).Value = 18
$ws.Cells.Item(55, 18).Value = "Hortaliza"
